$p = $ppt.ActivePresentation

# Remove the first slide, which is empty (no shapes) and is being
# dropped from the deck per the commit "removed first empty slide".
$p.Slides.Item(1).Delete()
